# re-run RU 1001; without crop
# This updates the simulation result values for the "Russia" (column L)
# scenario re-run, plus the small downstream recalculation effects on the
# "All" aggregate column (B) and a couple of floating point touch-ups in
# F6 and M8 caused by the re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2"  = 0.751312714286863
    "L2"  = 0.690938091049379
    "B3"  = 0.552311017972813
    "L3"  = 0.556331398643526
    "B4"  = 0.480678748076421
    "L4"  = 0.640825757278445
    "B5"  = 0.379373418323273
    "L5"  = 0.628339672766173
    "B6"  = 0.351492829218294
    "F6"  = 0.746324996119564
    "L6"  = -0.0300563026874272
    "B7"  = 0.35042789829748
    "L7"  = 0.672466186789019
    "B8"  = 0.344187928003983
    "L8"  = 0.353071433941653
    "M8"  = 0.627551511945664
    "B9"  = 0.341449989049043
    "L9"  = 0.370048816122248
    "B10" = 0.31420329928796
    "L10" = -0.0227328576067976
    "B11" = 0.0191486240382446
    "L11" = 0.0273835716769448
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
